$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# menambahkan 6 category baru
# Fill in Directory (C), By (D) and Status (E) columns for the newly
# completed categories. The write order below matches the order new
# unique strings were introduced in the shared-string table.

# Row 23: Gas
$ws.Range("C23").Value = "gas"

# Rows 25-29: Machinery, Power, Renovation, Contractor, Energy
$ws.Range("C25").Value = "machinery"
$ws.Range("C26").Value = "power"
$ws.Range("C27").Value = "renovation"
$ws.Range("C28").Value = "contractor"
$ws.Range("C29").Value = "energy"

# Row 30: Company (directory only, no owner/status yet)
$ws.Range("C30").Value = "company"

# Rows 32-39: Manufacturing, Building, Factory, Engineering, Business,
# Construction, Finance, Landing Page (directory only)
$ws.Range("C32").Value = "manufacturing"
$ws.Range("C33").Value = "building"
$ws.Range("C34").Value = "factory"
$ws.Range("C35").Value = "engineering"
$ws.Range("C36").Value = "business"
$ws.Range("C37").Value = "construction"
$ws.Range("C38").Value = "finance"
$ws.Range("C39").Value = "landing page"

# Row 23's "By" - a distinct "Krisna " (trailing space) entry
$ws.Range("D23").Value = "Krisna "

# Rows 40-46: Real Estate, Agency, Job, Consulting, Startup, Ecommerce,
# Industrial (directory only)
$ws.Range("C40").Value = "real estate"
$ws.Range("C41").Value = "agency"
$ws.Range("C42").Value = "job"
$ws.Range("C43").Value = "consulting"
$ws.Range("C44").Value = "startup"
$ws.Range("C45").Value = "ecommerce"
$ws.Range("C46").Value = "industrial"

# "By" / "Status" for rows 25-29 reuse the existing shared strings
$ws.Range("D25").Value = "Krisna"
$ws.Range("E23").Value = "DONE"
$ws.Range("E25").Value = "DONE"

$ws.Range("D26").Value = "Krisna"
$ws.Range("E26").Value = "DONE"

$ws.Range("D27").Value = "Krisna"
$ws.Range("E27").Value = "DONE"

$ws.Range("D28").Value = "Krisna"
$ws.Range("E28").Value = "DONE"

$ws.Range("D29").Value = "Krisna"
$ws.Range("E29").Value = "DONE"

# Restore view/selection state
$ws.Range("J27").Select()
